$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.778.14'
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").Value = '2.498.81'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'322.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = "'109.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("E7").Value = '  -0.62%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = "'0.554"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.79%  '
$ws.Range("D10").Value = "'40.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.80%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = "'0.124"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").Value = "'18.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.60%  '
$ws.Range("D14").Value = "'7.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("D15").Value = '2.892.60'
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").Value = '2.504.06'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = "'0.855"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = '47.664.59'
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("D19").Value = "'13.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.57%  '
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("D21").Value = "'2.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +13.83%  '
$ws.Range("D22").Value = '0.0₃0944'
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = "'247.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("D25").Value = "'2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = "'25.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("D28").Value = "'10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("D30").Value = "'0.140"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("D31").Value = "'35.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.04%  '
$ws.Range("D32").Value = "'49.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").Value = "'20.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").Value = "'5.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("D39").Value = "'2.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = "'22.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.22%  '
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("D43").Value = "'119.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("D44").Value = "'0.0299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("D45").Value = '2.002.98'
$ws.Range("E45").Value = '  +1.68%  '
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("E47").Value = '  -3.25%  '
$ws.Range("D48").Value = "'1.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").Value = "'9.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.13%  '
$ws.Range("D50").Value = "'5.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("D51").Value = "'56.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.93%  '
